# Edit the flow-chart slide (slide 1):
#  - Shape "직사각형 13" (id=14, 5th shape): "데이터 수집" -> "기상데이터 수집"
#    (typed "기상" in front of "데이터", splitting the original run into two
#    runs: "기상데이터" + " 수집").
#  - Shape "직사각형 29" (id=30, 7th shape): "데이터 삽입" -> two lines,
#    "데이터 " on the first line and "전처리" on the second line (the word
#    "삽입" was replaced by a paragraph break followed by "전처리").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "데이터 수집" -> "기상데이터" + " 수집" ------------------------------
$shp1 = $s.Shapes.Item(5)
$tr1 = $shp1.TextFrame.TextRange
$tr1.InsertBefore("기상데이터")
# The paragraph now holds two runs: "기상데이터" and the original "데이터 수집".
# Trim the leading "데이터" off the second (original) run so the visible
# text reads "기상데이터 수집".
$shp1.TextFrame.TextRange.Runs(2).Text = " 수집"

# --- "데이터 삽입" -> "데이터 " / "전처리" (two paragraphs) ---------------
$shp2 = $s.Shapes.Item(7)
$tr2 = $shp2.TextFrame.TextRange
$tr2.Text = "데이터 "
$tr2.InsertAfter("`r전처리")
